$wb = $excel.ActiveWorkbook

# Helper: apply a "top+bottom" thin border (matches new style index 2 / borderId 4)
function Set-TopBottomBorder($rng) {
    $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders.Item(8).Weight = 2      # xlThin
    $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders.Item(9).Weight = 2      # xlThin
}

# Helper: apply a "top+bottom+right" thin border (matches new style index 3 / borderId 5)
function Set-TopBottomRightBorder($rng) {
    $rng.Borders.Item(8).LineStyle = 1    # xlEdgeTop
    $rng.Borders.Item(8).Weight = 2       # xlThin
    $rng.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
    $rng.Borders.Item(9).Weight = 2       # xlThin
    $rng.Borders.Item(10).LineStyle = 1   # xlEdgeRight
    $rng.Borders.Item(10).Weight = 2      # xlThin
}

# ---------- Sheet 1: quality_comparison ----------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# ---------- Sheet 2: computational_comparison ----------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5 entirely
$ws2.Range("G5").ClearContents()
